$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has a header row (row 1) and one data row (row 2) for
# Mohammad Nabi's IPL innings. The diff adds a second data row (row 3) that
# duplicates row 2's record exactly (same match, same stats), extending the
# used range from A1:K2 to A1:K3.
#
# Force text formatting first so the numeric-looking values (runs, balls,
# 4s, 6s, strike rate) are written as text strings, matching how row 1/2
# already store every value as text (t="str" in the source XML / ignored
# "number stored as text" warning).
$ws.Range("A3:K3").NumberFormat = "@"

$ws.Range("A3").Value = " Abu Dhabi"
$ws.Range("B3").Value = " September 26 2020"
$ws.Range("C3").Value = "KKR won by 7 wickets (with 12 balls remaining)"
$ws.Range("D3").Value = "Sunrisers Hyderabad"
$ws.Range("E3").Value = "Kolkata Knight Riders"
$ws.Range("F3").Value = "Mohammad Nabi "
$ws.Range("G3").Value = "11"
$ws.Range("H3").Value = "8"
$ws.Range("I3").Value = "2"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "137.50"
